# COVID-19 Growth Projections workbook update
# "Update with some April 4,5 and 6 data."
#
# Sheet1 layout: A=date, B=projected r=1.189, C=projected r=1.3,
# D=CDC/JH reported cumulative cases, E=reported cumulative deaths,
# F=E/D (death rate), G=(E-Eprev)/(D-Dprev) (differential death rate).
#
# New data: row 25 = Apr-4, row 26 = Apr-5, row 27 = Apr-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 24 (Apr-3) already had D24/E24; add the derived ratios that were
#     missing for that day. ---
$ws.Range("F24").Formula = "=E24/D24"
$ws.Range("F24").NumberFormat = "0.00%"

$ws.Range("G24").Formula = "=(E24-E23)/(D24-D23)"
$ws.Range("G24").NumberFormat = "0.00%"

# --- Row 25 (Apr-4): only a confirmed-case count was available. ---
$ws.Range("D25").Value = 277205
$ws.Range("D25").NumberFormat = "#,##0"

# --- Row 26 (Apr-5): confirmed cases + deaths, plus derived ratios. ---
$ws.Range("D26").Value = 304826
$ws.Range("D26").NumberFormat = "#,##0"

$ws.Range("E26").Value = 7616
$ws.Range("E26").NumberFormat = "#,##0"

$ws.Range("F26:F27").Formula = "=E26/D26"
$ws.Range("F26:F27").NumberFormat = "0.00%"

$ws.Range("G26").Formula = "=(E26-E24)/(D26-D24)"
$ws.Range("G26").NumberFormat = "0.00%"

# --- Row 27 (Apr-6): confirmed cases + deaths, plus derived ratios. ---
$ws.Range("D27").Value = 330891
$ws.Range("D27").NumberFormat = "#,##0"

$ws.Range("E27").Value = 8910
$ws.Range("E27").NumberFormat = "#,##0"

$ws.Range("G27").Formula = "=(E27-E26)/(D27-D26)"
$ws.Range("G27").NumberFormat = "0.00%"

$excel.CalculateFull()
